# Apply the "task transformer" spreadsheet header changes:
#  - remove support for column "Version" (handled elsewhere / kept as-is here,
#    since the existing "Version" header cell AB1 is untouched in the diff)
#  - add support for columns "goal_version", "rule_name_id", "rule_version"
#    immediately after the existing "goal_name_id" column (AN), reusing the
#    three already-blank columns AO/AP/AQ that sat between "goal_name_id"
#    (AN) and the old "Parameter"/"Values" columns (AR/AS).
#  - insert two fresh blank columns where the old "Parameter"/"Values"
#    columns used to start, so that those two columns shift two places to
#    the right (AR,AS -> AT,AU), matching how real Excel pushed them over
#    when the new columns were added in front of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the trailing "Parameter" / "Values" columns two places right ---
# (old AR:AS -> new AT:AU). Inserting blank columns at AR:AS pushes any
# existing content (and its formatting) to the right automatically.
$ws.Range("AR:AS").Insert() | Out-Null

# --- Fill the (already existing, previously blank) AO/AP/AQ columns ---
# Row 1 headers: copy formatting from the neighbouring "goal_name_id"
# header (AN1) so the new headers render the same (bold, centered).
$ws.Range("AN1").Copy() | Out-Null
$ws.Range("AO1:AQ1").PasteSpecial(-4122) | Out-Null
$ws.Range("AO1").Value = "goal_version"
$ws.Range("AP1").Value = "rule_name_id"
$ws.Range("AQ1").Value = "rule_version"

# Row 2 data: copy formatting from the neighbouring "goal_name_id" value
# cell (AN2); only AP2 (rule_name_id) gets a value for this sample row,
# mirroring AN2's existing goal id value. AO2/AQ2 stay blank (but
# formatted) just like AG2 elsewhere on the sheet.
$ws.Range("AN2").Copy() | Out-Null
$ws.Range("AO2:AQ2").PasteSpecial(-4122) | Out-Null
$ws.Range("AP2").Value = $ws.Range("AN2").Value()

# AP (rule_name_id) ends up the same display width as AO (15 chars).
$ws.Range("AP1").ColumnWidth = 14.166666666666666

# --- Update the saved selection to reflect where the author ended up ---
$ws.Range("AR2").Select() | Out-Null
